$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 293.3
$ws.Range("I55").Value = 293.33334
$ws.Range("J55").Value = 293.25
$ws.Range("K55").Value = 293.33334
$ws.Range("L55").Value = 293.25
$ws.Range("M55").Value = -79.33334000000002
$ws.Range("N55").Value = -721.25
$ws.Range("H62").Value = 3769.9092
$ws.Range("I62").Value = 3360
$ws.Range("J62").Value = 4261.8
$ws.Range("K62").Value = 3360
$ws.Range("L62").Value = 4261.8
$ws.Range("M62").Value = -2736
$ws.Range("N62").Value = -5509.8
$ws.Range("H65").Value = 3769.9092
$ws.Range("I65").Value = 3360
$ws.Range("J65").Value = 4261.8
$ws.Range("K65").Value = 16800
$ws.Range("L65").Value = 21309
$ws.Range("M65").Value = -13680
$ws.Range("N65").Value = -27549
$ws.Range("H125").Value = 4919.2856
$ws.Range("I125").Value = 1865.1666
$ws.Range("J125").Value = 8991.444
$ws.Range("K125").Value = 16786.4994
$ws.Range("L125").Value = 80922.996
$ws.Range("M125").Value = -14326.4994
$ws.Range("N125").Value = -85842.996
$ws.Range("H137").Value = 217462.02
$ws.Range("I137").Value = 348924.78
$ws.Range("J137").Value = 5660.8887
$ws.Range("K137").Value = 1046774.34
$ws.Range("L137").Value = 16982.6661
$ws.Range("M137").Value = -1044224.34
$ws.Range("N137").Value = -22082.6661
$ws.Range("H141").Value = 2498.75
$ws.Range("I141").Value = 2498.75
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7496.25
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -2316.25
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20803.135
$ws.Range("I32").Value = 20182.562
$ws.Range("J32").Value = 28250
$ws.Range("K32").Value = 20182.562
$ws.Range("L32").Value = 28250
$ws.Range("M32").Value = -19895.562
$ws.Range("N32").Value = -28824
$ws.Range("H41").Value = 18818
$ws.Range("J41").Value = 25999.5
$ws.Range("L41").Value = 25999.5
$ws.Range("N41").Value = -26827.5
$ws.Range("H74").Value = 769.0769
$ws.Range("I74").Value = 560.3
$ws.Range("J74").Value = 1465
$ws.Range("K74").Value = 560.3
$ws.Range("L74").Value = 1465
$ws.Range("M74").Value = 313.7
$ws.Range("N74").Value = -3213
$ws.Range("H77").Value = 769.0769
$ws.Range("I77").Value = 560.3
$ws.Range("J77").Value = 1465
$ws.Range("K77").Value = 2801.5
$ws.Range("L77").Value = 7325
$ws.Range("M77").Value = 1566.5
$ws.Range("N77").Value = -16061
$ws.Range("H102").Value = 1609.75
$ws.Range("I102").Value = 1601.25
$ws.Range("K102").Value = 1601.25
$ws.Range("M102").Value = 20.75
$ws.Range("H122").Value = 4045.0908
$ws.Range("I122").Value = 2330.75
$ws.Range("J122").Value = 8616.666999999999
$ws.Range("K122").Value = 6992.25
$ws.Range("L122").Value = 25850.001
$ws.Range("M122").Value = -4542.25
$ws.Range("N122").Value = -30750.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4624.75
$ws.Range("I54").Value = 4624.75
$ws.Range("K54").Value = 4624.75
$ws.Range("M54").Value = -4140.75
$ws.Range("H99").Value = 3847.7778
$ws.Range("I99").Value = 3847.7778
$ws.Range("K99").Value = 3847.7778
$ws.Range("M99").Value = -2349.7778
$ws.Range("H105").Value = 2700.0476
$ws.Range("I105").Value = 2673.8948
$ws.Range("J105").Value = 2948.5
$ws.Range("K105").Value = 2673.8948
$ws.Range("L105").Value = 2948.5
$ws.Range("M105").Value = -926.8948
$ws.Range("N105").Value = -6442.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 2666.6667
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H16").Value = 2301.1667
$ws.Range("I16").Value = 1622.4445
$ws.Range("J16").Value = 4337.3335
$ws.Range("K16").Value = 1622.4445
$ws.Range("L16").Value = 4337.3335
$ws.Range("M16").Value = -1335.4445
$ws.Range("N16").Value = -4911.3335
$ws.Range("H31").Value = 33335880
$ws.Range("I31").Value = 35716372
$ws.Range("J31").Value = 8998.5
$ws.Range("K31").Value = 35716372
$ws.Range("L31").Value = 8998.5
$ws.Range("M31").Value = -35716077
$ws.Range("N31").Value = -9588.5
$ws.Range("H34").Value = 33335880
$ws.Range("I34").Value = 35716372
$ws.Range("J34").Value = 8998.5
$ws.Range("K34").Value = 35716372
$ws.Range("L34").Value = 8998.5
$ws.Range("M34").Value = -35716170
$ws.Range("N34").Value = -9402.5
$ws.Range("H58").Value = 3163.3572
$ws.Range("I58").Value = 2246
$ws.Range("J58").Value = 3851.375
$ws.Range("K58").Value = 2246
$ws.Range("L58").Value = 3851.375
$ws.Range("M58").Value = -2043
$ws.Range("N58").Value = -4257.375
$ws.Range("H113").Value = 2301.1667
$ws.Range("I113").Value = 1622.4445
$ws.Range("J113").Value = 4337.3335
$ws.Range("K113").Value = 1622.4445
$ws.Range("L113").Value = 4337.3335
$ws.Range("M113").Value = 547.5554999999999
$ws.Range("N113").Value = -8677.333500000001
$ws.Range("H132").Value = 63493564
$ws.Range("I132").Value = 63493564
$ws.Range("K132").Value = 190480692
$ws.Range("M132").Value = -190478162
$ws.Range("H136").Value = 3163.3572
$ws.Range("I136").Value = 2246
$ws.Range("J136").Value = 3851.375
$ws.Range("K136").Value = 6738
$ws.Range("L136").Value = 11554.125
$ws.Range("M136").Value = -4188
$ws.Range("N136").Value = -16654.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1243.4445
$ws.Range("I34").Value = 353.36365
$ws.Range("J34").Value = 2642.1428
$ws.Range("K34").Value = 1060.09095
$ws.Range("L34").Value = 7926.428400000001
$ws.Range("M34").Value = -976.09095
$ws.Range("N34").Value = -8094.428400000001
$ws.Range("H39").Value = 3908.818
$ws.Range("J39").Value = 3908.818
$ws.Range("L39").Value = 11726.454
$ws.Range("N39").Value = -12314.454
$ws.Range("H55").Value = 5334
$ws.Range("J55").Value = 5800
$ws.Range("L55").Value = 17400
$ws.Range("N55").Value = -17754
$ws.Range("H107").Value = 649.8
$ws.Range("I107").Value = 250
$ws.Range("K107").Value = 750
$ws.Range("M107").Value = 1170
$ws.Range("H122").Value = 1929.8
$ws.Range("J122").Value = 1824.25
$ws.Range("L122").Value = 16418.25
$ws.Range("N122").Value = -21318.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4500
$ws.Range("J70").Value = 4500
$ws.Range("L70").Value = 4500
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 4500
$ws.Range("J73").Value = 4500
$ws.Range("L73").Value = 4500
$ws.Range("N73").Value = -6372
$ws.Range("H122").Value = 275448.94
$ws.Range("I122").Value = 503827.9
$ws.Range("J122").Value = 6767.8237
$ws.Range("K122").Value = 1511483.7
$ws.Range("L122").Value = 20303.4711
$ws.Range("M122").Value = -1509033.7
$ws.Range("N122").Value = -25203.4711
$ws.Range("H132").Value = 56244.54
$ws.Range("I132").Value = 78750.08
$ws.Range("J132").Value = 3049.6365
$ws.Range("K132").Value = 236250.24
$ws.Range("L132").Value = 9148.9095
$ws.Range("M132").Value = -233720.24
$ws.Range("N132").Value = -14208.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2759.5251
$ws.Range("I132").Value = 2738.4492
$ws.Range("J132").Value = 2808
$ws.Range("K132").Value = 8215.347600000001
$ws.Range("L132").Value = 8424
$ws.Range("M132").Value = -5685.347600000001
$ws.Range("N132").Value = -13484

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1927.4048
$ws.Range("I132").Value = 922.36
$ws.Range("J132").Value = 3405.4119
$ws.Range("K132").Value = 2767.08
$ws.Range("L132").Value = 10216.2357
$ws.Range("M132").Value = -237.0799999999999
$ws.Range("N132").Value = -15276.2357
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
